$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Boolean outlier flag values for each row (F, G, H columns)
$values = @{
    2  = @($false, $false, $false)
    3  = @($false, $false, $false)
    4  = @($false, $false, $false)
    5  = @($false, $false, $false)
    6  = @($true,  $false, $false)
    7  = @($false, $false, $false)
    8  = @($false, $false, $false)
    9  = @($false, $false, $false)
    10 = @($false, $false, $false)
    11 = @($false, $false, $false)
    12 = @($false, $false, $false)
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 6).Value = $v[0]
    $ws.Cells.Item($row, 7).Value = $v[1]
    $ws.Cells.Item($row, 8).Value = $v[2]
}
